$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 963
$ws.Range("J3").Value = 8076
$ws.Range("K3").Value = 909
$ws.Range("I4").Value = 1784
$ws.Range("J4").Value = 1796
$ws.Range("K4").Value = 206
$ws.Range("K6").Value = 1266
$ws.Range("I7").Value = 26237
$ws.Range("K7").Value = 3397

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 100
$ws.Range("K15").Value = 21
$ws.Range("K19").Value = 91
$ws.Range("J20").Value = 634
$ws.Range("K20").Value = 83
$ws.Range("K23").Value = 30
$ws.Range("K24").Value = 11
$ws.Range("K27").Value = 43
$ws.Range("K29").Value = 161
$ws.Range("K30").Value = 9
$ws.Range("K31").Value = 35
$ws.Range("K33").Value = 141
$ws.Range("J36").Value = 403
$ws.Range("K37").Value = 106
$ws.Range("K41").Value = 35
$ws.Range("K42").Value = 109
$ws.Range("K44").Value = 30
$ws.Range("J49").Value = 177
$ws.Range("K49").Value = 22
$ws.Range("K50").Value = 23
$ws.Range("K52").Value = 82
$ws.Range("J53").Value = 440
$ws.Range("K53").Value = 40
$ws.Range("K54").Value = 60
$ws.Range("I63").Value = 193
$ws.Range("J63").Value = 89
$ws.Range("K63").Value = 15
$ws.Range("J65").Value = 727
$ws.Range("K65").Value = 94
$ws.Range("K67").Value = 143
$ws.Range("K69").Value = 12
$ws.Range("J77").Value = 205
$ws.Range("K77").Value = 23
$ws.Range("K78").Value = 50
$ws.Range("K79").Value = 96
$ws.Range("K85").Value = 172
$ws.Range("K86").Value = 24
$ws.Range("K89").Value = 54
$ws.Range("K90").Value = 29
$ws.Range("K91").Value = 36
$ws.Range("K93").Value = 12
$ws.Range("K94").Value = 41
$ws.Range("K96").Value = 56
$ws.Range("K97").Value = 26
$ws.Range("I101").Value = 26237
$ws.Range("K101").Value = 3397

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 42
$ws.Range("K3").Value = 25
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 68
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 20
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 79
$ws.Range("K3").Value = 9
$ws.Range("J7").Value = 440
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 36
$ws.Range("K3").Value = 53
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 21
$ws.Range("K3").Value = 34
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 22
$ws.Range("K3").Value = 21
$ws.Range("J4").Value = 29
$ws.Range("J7").Value = 727
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 40
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 4
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 177
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 21
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 40
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 22
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 4
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 34
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 19
$ws.Range("J4").Value = 52
$ws.Range("J7").Value = 634
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 403

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 9
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 205
$ws.Range("K7").Value = 23
